$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trade row (row 4), mirroring the existing row 3 formatting.
# Copy A3's format (date number format) onto A4 before setting its value.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 42641.546655092592

$ws.Range("B4").Value = $false

$ws.Range("C4").Value = 9942.5300000000007
$ws.Range("D4").Value = 9948
$ws.Range("E4").Value = 18.670000000000002
$ws.Range("F4").Value = 18.690000000000001

# Copy G3's format (boolean column style) onto G4 before setting its value.
$ws.Range("G3").Copy()
$ws.Range("G4").PasteSpecial(-4122)
$ws.Range("G4").Value = $true

$ws.Range("H4").Value = 0.11
$ws.Range("I4").Value = $false
